$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 191; this shifts the existing rows 191-295
# down to 192-296 and extends the used range to A1:T296.
$ws.Rows.Item(191).Insert()

# Populate the newly inserted row 191 with the new price record.
$ws.Range("A191").Value = 4
$ws.Range("B191").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C191").Value = "Los Lagos"
$ws.Range("D191").Value = (Get-Date -Year 2022 -Month 9 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E191").Value = 10
$ws.Range("F191").Value = "Fruta"
$ws.Range("G191").Value = 100108
$ws.Range("H191").Value = "Tropicales y subtropicales"
$ws.Range("I191").Value = 100108005
$ws.Range("J191").Value = "Piña"
$ws.Range("K191").Value = "Caramelo"
$ws.Range("L191").Value = "Primera"
$ws.Range("M191").Value = 150
$ws.Range("N191").Value = 23000
$ws.Range("O191").Value = 23000
$ws.Range("P191").Value = 23000
$ws.Range("Q191").Value = "`$/caja 12 unidades"
$ws.Range("R191").Value = "Ecuador"
$ws.Range("S191").Value = 1917
$ws.Range("T191").Value = 12
